$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 307, pushing all
# existing data (old rows 307..355) down to 309..357.
$ws.Rows.Item(307).Insert()
$ws.Rows.Item(307).Insert()

# --- New row 307: Uva / Red Globe, Región de O'Higgins ---
$ws.Range("A307").Value = 4
$ws.Range("B307").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C307").Value = "Los Lagos"
$ws.Range("D307").Value = 45034
$ws.Range("E307").Value = 10
$ws.Range("F307").Value = "Fruta"
$ws.Range("G307").Value = 100109
$ws.Range("H307").Value = "Uva"
$ws.Range("I307").Value = 100109001
$ws.Range("J307").Value = "Uva"
$ws.Range("K307").Value = "Red Globe"
$ws.Range("L307").Value = "Primera"
$ws.Range("M307").Value = 400
$ws.Range("N307").Value = 14000
$ws.Range("O307").Value = 15000
$ws.Range("P307").Value = 14500
$ws.Range("Q307").Value = "`$/caja 18 kilos"
$ws.Range("R307").Value = "Región de O'Higgins"
$ws.Range("S307").Value = 806
$ws.Range("T307").Value = 18

# --- New row 308: Uva / Rosada pastilla, Región de O'Higgins ---
$ws.Range("A308").Value = 4
$ws.Range("B308").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C308").Value = "Los Lagos"
$ws.Range("D308").Value = 45034
$ws.Range("E308").Value = 10
$ws.Range("F308").Value = "Fruta"
$ws.Range("G308").Value = 100109
$ws.Range("H308").Value = "Uva"
$ws.Range("I308").Value = 100109001
$ws.Range("J308").Value = "Uva"
$ws.Range("K308").Value = "Rosada pastilla"
$ws.Range("L308").Value = "Primera"
$ws.Range("M308").Value = 400
$ws.Range("N308").Value = 14000
$ws.Range("O308").Value = 15000
$ws.Range("P308").Value = 14500
$ws.Range("Q308").Value = "`$/bandeja 10 kilos"
$ws.Range("R308").Value = "Región de O'Higgins"
$ws.Range("S308").Value = 1450
$ws.Range("T308").Value = 10
